$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TestCase1_HospitalFilter")
$ws.Cells.Item(127, 1).Value = 'Opened Practo homepage.'
$ws.Cells.Item(128, 1).Value = 'Searching for hospitals in: Bangalore'
$ws.Cells.Item(129, 1).Value = 'Added hospital: Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws.Cells.Item(130, 1).Value = 'Added hospital: Manipal Hospitals'
$ws.Cells.Item(131, 1).Value = 'Added hospital: Koshys Hospital'
$ws.Cells.Item(132, 1).Value = 'Added hospital: Motherhood Hospital'
$ws.Cells.Item(133, 1).Value = 'Added hospital: Motherhood Hospital'
$ws.Cells.Item(134, 1).Value = 'Added hospital: Trilife Hospital'
$ws.Cells.Item(135, 1).Value = 'Added hospital: Apollo Cradle & Children’s Hospital'
$ws.Cells.Item(136, 1).Value = 'Printing and Writing'
$ws.Cells.Item(137, 1).Value = 'Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws.Cells.Item(138, 1).Value = 'Manipal Hospitals'
$ws.Cells.Item(139, 1).Value = 'Koshys Hospital'
$ws.Cells.Item(140, 1).Value = 'Motherhood Hospital'
$ws.Cells.Item(141, 1).Value = 'Motherhood Hospital'
$ws.Cells.Item(142, 1).Value = 'Trilife Hospital'
$ws.Cells.Item(143, 1).Value = 'Apollo Cradle & Children’s Hospital'
$ws.Cells.Item(144, 1).Value = 'Opened Practo homepage.'
$ws.Cells.Item(145, 1).Value = 'Searching for hospitals in: Bangalore'
$ws.Cells.Item(146, 1).Value = 'Added hospital: Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws.Cells.Item(147, 1).Value = 'Added hospital: Manipal Hospitals'
$ws.Cells.Item(148, 1).Value = 'Added hospital: Koshys Hospital'
$ws.Cells.Item(149, 1).Value = 'Added hospital: Motherhood Hospital'
$ws.Cells.Item(150, 1).Value = 'Added hospital: Motherhood Hospital'
$ws.Cells.Item(151, 1).Value = 'Added hospital: Trilife Hospital'
$ws.Cells.Item(152, 1).Value = 'Added hospital: Apollo Cradle & Children’s Hospital'
$ws.Cells.Item(153, 1).Value = 'Printing and Writing'
$ws.Cells.Item(154, 1).Value = 'Manipal Hospital Varthur Road (formerly Columbia Asia Hospital)'
$ws.Cells.Item(155, 1).Value = 'Manipal Hospitals'
$ws.Cells.Item(156, 1).Value = 'Koshys Hospital'
$ws.Cells.Item(157, 1).Value = 'Motherhood Hospital'
$ws.Cells.Item(158, 1).Value = 'Motherhood Hospital'
$ws.Cells.Item(159, 1).Value = 'Trilife Hospital'
$ws.Cells.Item(160, 1).Value = 'Apollo Cradle & Children’s Hospital'

$ws = $wb.Worksheets.Item("TestCase0_MaxRatingFinder")
$ws.Cells.Item(49, 1).Value = 'Searching for hospitals in: Bangalore'
$ws.Cells.Item(50, 1).Value = 'Error during search: Expected condition failed: waiting for element found by By.xpath: //div[contains(@class,''c-omni-suggestion-item'')]//div[contains(text(),''Bangalore'')] to be clickable, but the element was not found: org.openqa.selenium.NoSuchElementException: no such element: Unable to locate element: {"method":"xpath","selector":"//div[contains(@class,''c-omni-suggestion-item'')]//div[contains(text(),''Bangalore'')]"}.
(tried for 15 seconds with 500 milliseconds interval)
Build info: version: ''4.40.0'', revision: ''b3333f1''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''21''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 142.0.7444.176, chrome: {chromedriverVersion: 142.0.7444.175 (302067f14a4..., userDataDir: C:\Users\2457382\AppData\Lo...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:52484}, goog:processID: 14520, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:52484/devtoo..., se:cdpVersion: 142.0.7444.176, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 9d4310f28ca2e78faf41ff59fa9589c0'
$ws.Cells.Item(51, 1).Value = 'Checking for rating ≥ 4.50: Bangalore'
$ws.Cells.Item(52, 1).Value = 'Searching for hospitals in: Bangalore'
$ws.Cells.Item(53, 1).Value = 'Checking for rating ≥ 4.50: Bangalore'
$ws.Cells.Item(54, 1).Value = 'Rating captured: 4.50'
$ws.Cells.Item(55, 1).Value = 'PASS — Rating ≥ 4.50 (actual: 4.50)'
$ws.Cells.Item(56, 1).Value = 'Result: Max Rating meets threshold (≥ 4.50): actual 4.50'

$ws = $wb.Worksheets.Item("TestCase2_TopCities")
$ws.Cells.Item(73, 1).Value = 'Top Cities:'
$ws.Cells.Item(74, 1).Value = 'Bangalore'
$ws.Cells.Item(75, 1).Value = 'Delhi'
$ws.Cells.Item(76, 1).Value = 'Mumbai'
$ws.Cells.Item(77, 1).Value = 'Chennai'
$ws.Cells.Item(78, 1).Value = 'Hyderabad'
$ws.Cells.Item(79, 1).Value = 'Kolkata'
$ws.Cells.Item(80, 1).Value = 'Pune'
$ws.Cells.Item(81, 1).Value = 'Ahmedabad'
$ws.Cells.Item(82, 1).Value = 'Top Cities:'
$ws.Cells.Item(83, 1).Value = 'Bangalore'
$ws.Cells.Item(84, 1).Value = 'Delhi'
$ws.Cells.Item(85, 1).Value = 'Mumbai'
$ws.Cells.Item(86, 1).Value = 'Chennai'
$ws.Cells.Item(87, 1).Value = 'Hyderabad'
$ws.Cells.Item(88, 1).Value = 'Kolkata'
$ws.Cells.Item(89, 1).Value = 'Pune'
$ws.Cells.Item(90, 1).Value = 'Ahmedabad'
$ws.Cells.Item(91, 1).Value = 'Top Cities:'
$ws.Cells.Item(92, 1).Value = 'Bangalore'
$ws.Cells.Item(93, 1).Value = 'Delhi'
$ws.Cells.Item(94, 1).Value = 'Mumbai'
$ws.Cells.Item(95, 1).Value = 'Chennai'
$ws.Cells.Item(96, 1).Value = 'Hyderabad'
$ws.Cells.Item(97, 1).Value = 'Kolkata'
$ws.Cells.Item(98, 1).Value = 'Pune'
$ws.Cells.Item(99, 1).Value = 'Ahmedabad'

$ws = $wb.Worksheets.Item("TestCase3_InvalidForm")
$ws.Cells.Item(63, 1).Value = 'Forced click on Schedule button.'
$ws.Cells.Item(64, 1).Value = 'No error messages found.'
$ws.Cells.Item(65, 1).Value = 'Checking for invalid fields...'
$ws.Cells.Item(66, 1).Value = 'Empty Name'
$ws.Cells.Item(67, 1).Value = 'Empty Organization Name'
$ws.Cells.Item(68, 1).Value = 'Invalid Contact Number'
$ws.Cells.Item(69, 1).Value = 'Invalid Email ID'
$ws.Cells.Item(70, 1).Value = 'Forced click on Schedule button.'
$ws.Cells.Item(71, 1).Value = 'No error messages found.'
$ws.Cells.Item(72, 1).Value = 'Checking for invalid fields...'
$ws.Cells.Item(73, 1).Value = 'Empty Name'
$ws.Cells.Item(74, 1).Value = 'Empty Organization Name'
$ws.Cells.Item(75, 1).Value = 'Invalid Contact Number'
$ws.Cells.Item(76, 1).Value = 'Invalid Email ID'
$ws.Cells.Item(77, 1).Value = 'Forced click on Schedule button.'
$ws.Cells.Item(78, 1).Value = 'No error messages found.'
$ws.Cells.Item(79, 1).Value = 'Checking for invalid fields...'
$ws.Cells.Item(80, 1).Value = 'Empty Name'
$ws.Cells.Item(81, 1).Value = 'Empty Organization Name'
$ws.Cells.Item(82, 1).Value = 'Invalid Contact Number'
$ws.Cells.Item(83, 1).Value = 'Invalid Email ID'

